$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure new row 73 date cell uses same number format/style as other date cells in column D
$ws.Cells.Item(73, 4).NumberFormat = $ws.Cells.Item(72, 4).NumberFormat

# Row 31
$ws.Cells.Item(31, 1).Value2 = 1
$ws.Cells.Item(31, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value2 = 44895
$ws.Cells.Item(31, 5).Value2 = 15
$ws.Cells.Item(31, 6).Value2 = 100112031
$ws.Cells.Item(31, 7).Value2 = "Poroto verde"
$ws.Cells.Item(31, 8).Value2 = "Sin especificar"
$ws.Cells.Item(31, 9).Value2 = "Primera"
$ws.Cells.Item(31, 10).Value2 = 150
$ws.Cells.Item(31, 11).Value2 = 1300
$ws.Cells.Item(31, 12).Value2 = 1300
$ws.Cells.Item(31, 13).Value2 = 1300
$ws.Cells.Item(31, 14).Value2 = "$/kilo"
$ws.Cells.Item(31, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value2 = 1300
$ws.Cells.Item(31, 17).Value2 = 1
$ws.Cells.Item(31, 18).Value2 = "Hortaliza"

# Row 32
$ws.Cells.Item(32, 1).Value2 = 1
$ws.Cells.Item(32, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(32, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(32, 4).Value2 = 44873
$ws.Cells.Item(32, 5).Value2 = 15
$ws.Cells.Item(32, 6).Value2 = 100112031
$ws.Cells.Item(32, 7).Value2 = "Poroto verde"
$ws.Cells.Item(32, 8).Value2 = "Sin especificar"
$ws.Cells.Item(32, 9).Value2 = "Primera"
$ws.Cells.Item(32, 10).Value2 = 1200
$ws.Cells.Item(32, 11).Value2 = 800
$ws.Cells.Item(32, 12).Value2 = 900
$ws.Cells.Item(32, 13).Value2 = 850
$ws.Cells.Item(32, 14).Value2 = "$/kilo"
$ws.Cells.Item(32, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(32, 16).Value2 = 850
$ws.Cells.Item(32, 17).Value2 = 1
$ws.Cells.Item(32, 18).Value2 = "Hortaliza"

# Row 33
$ws.Cells.Item(33, 1).Value2 = 1
$ws.Cells.Item(33, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(33, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(33, 4).Value2 = 44312
$ws.Cells.Item(33, 5).Value2 = 15
$ws.Cells.Item(33, 6).Value2 = 100112031
$ws.Cells.Item(33, 7).Value2 = "Poroto verde"
$ws.Cells.Item(33, 8).Value2 = "Sin especificar"
$ws.Cells.Item(33, 9).Value2 = "Primera"
$ws.Cells.Item(33, 10).Value2 = 1700
$ws.Cells.Item(33, 11).Value2 = 1300
$ws.Cells.Item(33, 12).Value2 = 1400
$ws.Cells.Item(33, 13).Value2 = 1350
$ws.Cells.Item(33, 14).Value2 = "$/kilo"
$ws.Cells.Item(33, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(33, 16).Value2 = 1350
$ws.Cells.Item(33, 17).Value2 = 1
$ws.Cells.Item(33, 18).Value2 = "Hortaliza"

# Row 34
$ws.Cells.Item(34, 1).Value2 = 1
$ws.Cells.Item(34, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(34, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(34, 4).Value2 = 44760
$ws.Cells.Item(34, 5).Value2 = 15
$ws.Cells.Item(34, 6).Value2 = 100112031
$ws.Cells.Item(34, 7).Value2 = "Poroto verde"
$ws.Cells.Item(34, 8).Value2 = "Magnum"
$ws.Cells.Item(34, 9).Value2 = "Primera"
$ws.Cells.Item(34, 10).Value2 = 300
$ws.Cells.Item(34, 11).Value2 = 25000
$ws.Cells.Item(34, 12).Value2 = 26000
$ws.Cells.Item(34, 13).Value2 = 25500
$ws.Cells.Item(34, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(34, 15).Value2 = "Perú"
$ws.Cells.Item(34, 16).Value2 = 1020
$ws.Cells.Item(34, 17).Value2 = 25
$ws.Cells.Item(34, 18).Value2 = "Hortaliza"

# Row 35
$ws.Cells.Item(35, 1).Value2 = 1
$ws.Cells.Item(35, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(35, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(35, 4).Value2 = 44258
$ws.Cells.Item(35, 5).Value2 = 15
$ws.Cells.Item(35, 6).Value2 = 100112031
$ws.Cells.Item(35, 7).Value2 = "Poroto verde"
$ws.Cells.Item(35, 8).Value2 = "Sin especificar"
$ws.Cells.Item(35, 9).Value2 = "Primera"
$ws.Cells.Item(35, 10).Value2 = 1600
$ws.Cells.Item(35, 11).Value2 = 2300
$ws.Cells.Item(35, 12).Value2 = 2500
$ws.Cells.Item(35, 13).Value2 = 2400
$ws.Cells.Item(35, 14).Value2 = "$/kilo"
$ws.Cells.Item(35, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(35, 16).Value2 = 2400
$ws.Cells.Item(35, 17).Value2 = 1
$ws.Cells.Item(35, 18).Value2 = "Hortaliza"

# Row 36
$ws.Cells.Item(36, 1).Value2 = 1
$ws.Cells.Item(36, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(36, 4).Value2 = 44567
$ws.Cells.Item(36, 5).Value2 = 15
$ws.Cells.Item(36, 6).Value2 = 100112031
$ws.Cells.Item(36, 7).Value2 = "Poroto verde"
$ws.Cells.Item(36, 8).Value2 = "Sin especificar"
$ws.Cells.Item(36, 9).Value2 = "Primera"
$ws.Cells.Item(36, 10).Value2 = 1200
$ws.Cells.Item(36, 11).Value2 = 400
$ws.Cells.Item(36, 12).Value2 = 500
$ws.Cells.Item(36, 13).Value2 = 450
$ws.Cells.Item(36, 14).Value2 = "$/kilo"
$ws.Cells.Item(36, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(36, 16).Value2 = 450
$ws.Cells.Item(36, 17).Value2 = 1
$ws.Cells.Item(36, 18).Value2 = "Hortaliza"

# Row 37
$ws.Cells.Item(37, 1).Value2 = 1
$ws.Cells.Item(37, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(37, 4).Value2 = 44330
$ws.Cells.Item(37, 5).Value2 = 15
$ws.Cells.Item(37, 6).Value2 = 100112031
$ws.Cells.Item(37, 7).Value2 = "Poroto verde"
$ws.Cells.Item(37, 8).Value2 = "Magnum"
$ws.Cells.Item(37, 9).Value2 = "Primera"
$ws.Cells.Item(37, 10).Value2 = 200
$ws.Cells.Item(37, 11).Value2 = 24000
$ws.Cells.Item(37, 12).Value2 = 25000
$ws.Cells.Item(37, 13).Value2 = 24500
$ws.Cells.Item(37, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(37, 15).Value2 = "Perú"
$ws.Cells.Item(37, 16).Value2 = 980
$ws.Cells.Item(37, 17).Value2 = 25
$ws.Cells.Item(37, 18).Value2 = "Hortaliza"

# Row 38
$ws.Cells.Item(38, 1).Value2 = 1
$ws.Cells.Item(38, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(38, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(38, 4).Value2 = 44348
$ws.Cells.Item(38, 5).Value2 = 15
$ws.Cells.Item(38, 6).Value2 = 100112031
$ws.Cells.Item(38, 7).Value2 = "Poroto verde"
$ws.Cells.Item(38, 8).Value2 = "Sin especificar"
$ws.Cells.Item(38, 9).Value2 = "Primera"
$ws.Cells.Item(38, 10).Value2 = 1300
$ws.Cells.Item(38, 11).Value2 = 1700
$ws.Cells.Item(38, 12).Value2 = 1800
$ws.Cells.Item(38, 13).Value2 = 1750
$ws.Cells.Item(38, 14).Value2 = "$/kilo"
$ws.Cells.Item(38, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(38, 16).Value2 = 1750
$ws.Cells.Item(38, 17).Value2 = 1
$ws.Cells.Item(38, 18).Value2 = "Hortaliza"

# Row 39
$ws.Cells.Item(39, 1).Value2 = 1
$ws.Cells.Item(39, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(39, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(39, 4).Value2 = 44399
$ws.Cells.Item(39, 5).Value2 = 15
$ws.Cells.Item(39, 6).Value2 = 100112031
$ws.Cells.Item(39, 7).Value2 = "Poroto verde"
$ws.Cells.Item(39, 8).Value2 = "Magnum"
$ws.Cells.Item(39, 9).Value2 = "Primera"
$ws.Cells.Item(39, 10).Value2 = 1400
$ws.Cells.Item(39, 11).Value2 = 1300
$ws.Cells.Item(39, 12).Value2 = 1400
$ws.Cells.Item(39, 13).Value2 = 1350
$ws.Cells.Item(39, 14).Value2 = "$/kilo"
$ws.Cells.Item(39, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 16).Value2 = 1350
$ws.Cells.Item(39, 17).Value2 = 1
$ws.Cells.Item(39, 18).Value2 = "Hortaliza"

# Row 40
$ws.Cells.Item(40, 1).Value2 = 1
$ws.Cells.Item(40, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(40, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(40, 4).Value2 = 44473
$ws.Cells.Item(40, 5).Value2 = 15
$ws.Cells.Item(40, 6).Value2 = 100112031
$ws.Cells.Item(40, 7).Value2 = "Poroto verde"
$ws.Cells.Item(40, 8).Value2 = "Sin especificar"
$ws.Cells.Item(40, 9).Value2 = "Primera"
$ws.Cells.Item(40, 10).Value2 = 1000
$ws.Cells.Item(40, 11).Value2 = 1600
$ws.Cells.Item(40, 12).Value2 = 1700
$ws.Cells.Item(40, 13).Value2 = 1650
$ws.Cells.Item(40, 14).Value2 = "$/kilo"
$ws.Cells.Item(40, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(40, 16).Value2 = 1650
$ws.Cells.Item(40, 17).Value2 = 1
$ws.Cells.Item(40, 18).Value2 = "Hortaliza"

# Row 41
$ws.Cells.Item(41, 1).Value2 = 1
$ws.Cells.Item(41, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(41, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(41, 4).Value2 = 44357
$ws.Cells.Item(41, 5).Value2 = 15
$ws.Cells.Item(41, 6).Value2 = 100112031
$ws.Cells.Item(41, 7).Value2 = "Poroto verde"
$ws.Cells.Item(41, 8).Value2 = "Magnum"
$ws.Cells.Item(41, 9).Value2 = "Primera"
$ws.Cells.Item(41, 10).Value2 = 300
$ws.Cells.Item(41, 11).Value2 = 22000
$ws.Cells.Item(41, 12).Value2 = 23000
$ws.Cells.Item(41, 13).Value2 = 22500
$ws.Cells.Item(41, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(41, 15).Value2 = "Perú"
$ws.Cells.Item(41, 16).Value2 = 900
$ws.Cells.Item(41, 17).Value2 = 25
$ws.Cells.Item(41, 18).Value2 = "Hortaliza"

# Row 42
$ws.Cells.Item(42, 1).Value2 = 1
$ws.Cells.Item(42, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(42, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(42, 4).Value2 = 44859
$ws.Cells.Item(42, 5).Value2 = 15
$ws.Cells.Item(42, 6).Value2 = 100112031
$ws.Cells.Item(42, 7).Value2 = "Poroto verde"
$ws.Cells.Item(42, 8).Value2 = "Sin especificar"
$ws.Cells.Item(42, 9).Value2 = "Primera"
$ws.Cells.Item(42, 10).Value2 = 1800
$ws.Cells.Item(42, 11).Value2 = 1400
$ws.Cells.Item(42, 12).Value2 = 1500
$ws.Cells.Item(42, 13).Value2 = 1444
$ws.Cells.Item(42, 14).Value2 = "$/kilo"
$ws.Cells.Item(42, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(42, 16).Value2 = 1444
$ws.Cells.Item(42, 17).Value2 = 1
$ws.Cells.Item(42, 18).Value2 = "Hortaliza"

# Row 43
$ws.Cells.Item(43, 1).Value2 = 1
$ws.Cells.Item(43, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(43, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(43, 4).Value2 = 44342
$ws.Cells.Item(43, 5).Value2 = 15
$ws.Cells.Item(43, 6).Value2 = 100112031
$ws.Cells.Item(43, 7).Value2 = "Poroto verde"
$ws.Cells.Item(43, 8).Value2 = "Sin especificar"
$ws.Cells.Item(43, 9).Value2 = "Primera"
$ws.Cells.Item(43, 10).Value2 = 1500
$ws.Cells.Item(43, 11).Value2 = 1300
$ws.Cells.Item(43, 12).Value2 = 1400
$ws.Cells.Item(43, 13).Value2 = 1350
$ws.Cells.Item(43, 14).Value2 = "$/kilo"
$ws.Cells.Item(43, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(43, 16).Value2 = 1350
$ws.Cells.Item(43, 17).Value2 = 1
$ws.Cells.Item(43, 18).Value2 = "Hortaliza"

# Row 44
$ws.Cells.Item(44, 1).Value2 = 1
$ws.Cells.Item(44, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(44, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(44, 4).Value2 = 44356
$ws.Cells.Item(44, 5).Value2 = 15
$ws.Cells.Item(44, 6).Value2 = 100112031
$ws.Cells.Item(44, 7).Value2 = "Poroto verde"
$ws.Cells.Item(44, 8).Value2 = "Sin especificar"
$ws.Cells.Item(44, 9).Value2 = "Primera"
$ws.Cells.Item(44, 10).Value2 = 1100
$ws.Cells.Item(44, 11).Value2 = 1400
$ws.Cells.Item(44, 12).Value2 = 1500
$ws.Cells.Item(44, 13).Value2 = 1450
$ws.Cells.Item(44, 14).Value2 = "$/kilo"
$ws.Cells.Item(44, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(44, 16).Value2 = 1450
$ws.Cells.Item(44, 17).Value2 = 1
$ws.Cells.Item(44, 18).Value2 = "Hortaliza"

# Row 45
$ws.Cells.Item(45, 1).Value2 = 1
$ws.Cells.Item(45, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(45, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(45, 4).Value2 = 44879
$ws.Cells.Item(45, 5).Value2 = 15
$ws.Cells.Item(45, 6).Value2 = 100112031
$ws.Cells.Item(45, 7).Value2 = "Poroto verde"
$ws.Cells.Item(45, 8).Value2 = "Sin especificar"
$ws.Cells.Item(45, 9).Value2 = "Primera"
$ws.Cells.Item(45, 10).Value2 = 1300
$ws.Cells.Item(45, 11).Value2 = 1300
$ws.Cells.Item(45, 12).Value2 = 1400
$ws.Cells.Item(45, 13).Value2 = 1350
$ws.Cells.Item(45, 14).Value2 = "$/kilo"
$ws.Cells.Item(45, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 16).Value2 = 1350
$ws.Cells.Item(45, 17).Value2 = 1
$ws.Cells.Item(45, 18).Value2 = "Hortaliza"

# Row 46
$ws.Cells.Item(46, 1).Value2 = 1
$ws.Cells.Item(46, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value2 = 44340
$ws.Cells.Item(46, 5).Value2 = 15
$ws.Cells.Item(46, 6).Value2 = 100112031
$ws.Cells.Item(46, 7).Value2 = "Poroto verde"
$ws.Cells.Item(46, 8).Value2 = "Magnum"
$ws.Cells.Item(46, 9).Value2 = "Primera"
$ws.Cells.Item(46, 10).Value2 = 600
$ws.Cells.Item(46, 11).Value2 = 24000
$ws.Cells.Item(46, 12).Value2 = 25000
$ws.Cells.Item(46, 13).Value2 = 24500
$ws.Cells.Item(46, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(46, 15).Value2 = "Perú"
$ws.Cells.Item(46, 16).Value2 = 980
$ws.Cells.Item(46, 17).Value2 = 25
$ws.Cells.Item(46, 18).Value2 = "Hortaliza"

# Row 47
$ws.Cells.Item(47, 1).Value2 = 1
$ws.Cells.Item(47, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(47, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(47, 4).Value2 = 44596
$ws.Cells.Item(47, 5).Value2 = 15
$ws.Cells.Item(47, 6).Value2 = 100112031
$ws.Cells.Item(47, 7).Value2 = "Poroto verde"
$ws.Cells.Item(47, 8).Value2 = "Sin especificar"
$ws.Cells.Item(47, 9).Value2 = "Primera"
$ws.Cells.Item(47, 10).Value2 = 1300
$ws.Cells.Item(47, 11).Value2 = 1400
$ws.Cells.Item(47, 12).Value2 = 1500
$ws.Cells.Item(47, 13).Value2 = 1450
$ws.Cells.Item(47, 14).Value2 = "$/kilo"
$ws.Cells.Item(47, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(47, 16).Value2 = 1450
$ws.Cells.Item(47, 17).Value2 = 1
$ws.Cells.Item(47, 18).Value2 = "Hortaliza"

# Row 48
$ws.Cells.Item(48, 1).Value2 = 1
$ws.Cells.Item(48, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(48, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(48, 4).Value2 = 44792
$ws.Cells.Item(48, 5).Value2 = 15
$ws.Cells.Item(48, 6).Value2 = 100112031
$ws.Cells.Item(48, 7).Value2 = "Poroto verde"
$ws.Cells.Item(48, 8).Value2 = "Magnum"
$ws.Cells.Item(48, 9).Value2 = "Primera"
$ws.Cells.Item(48, 10).Value2 = 160
$ws.Cells.Item(48, 11).Value2 = 24000
$ws.Cells.Item(48, 12).Value2 = 25000
$ws.Cells.Item(48, 13).Value2 = 24500
$ws.Cells.Item(48, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value2 = "Perú"
$ws.Cells.Item(48, 16).Value2 = 980
$ws.Cells.Item(48, 17).Value2 = 25
$ws.Cells.Item(48, 18).Value2 = "Hortaliza"

# Row 49
$ws.Cells.Item(49, 1).Value2 = 1
$ws.Cells.Item(49, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(49, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(49, 4).Value2 = 44313
$ws.Cells.Item(49, 5).Value2 = 15
$ws.Cells.Item(49, 6).Value2 = 100112031
$ws.Cells.Item(49, 7).Value2 = "Poroto verde"
$ws.Cells.Item(49, 8).Value2 = "Sin especificar"
$ws.Cells.Item(49, 9).Value2 = "Primera"
$ws.Cells.Item(49, 10).Value2 = 1900
$ws.Cells.Item(49, 11).Value2 = 1000
$ws.Cells.Item(49, 12).Value2 = 1200
$ws.Cells.Item(49, 13).Value2 = 1100
$ws.Cells.Item(49, 14).Value2 = "$/kilo"
$ws.Cells.Item(49, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(49, 16).Value2 = 1100
$ws.Cells.Item(49, 17).Value2 = 1
$ws.Cells.Item(49, 18).Value2 = "Hortaliza"

# Row 50
$ws.Cells.Item(50, 1).Value2 = 1
$ws.Cells.Item(50, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(50, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(50, 4).Value2 = 44537
$ws.Cells.Item(50, 5).Value2 = 15
$ws.Cells.Item(50, 6).Value2 = 100112031
$ws.Cells.Item(50, 7).Value2 = "Poroto verde"
$ws.Cells.Item(50, 8).Value2 = "Sin especificar"
$ws.Cells.Item(50, 9).Value2 = "Primera"
$ws.Cells.Item(50, 10).Value2 = 1700
$ws.Cells.Item(50, 11).Value2 = 500
$ws.Cells.Item(50, 12).Value2 = 600
$ws.Cells.Item(50, 13).Value2 = 550
$ws.Cells.Item(50, 14).Value2 = "$/kilo"
$ws.Cells.Item(50, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(50, 16).Value2 = 550
$ws.Cells.Item(50, 17).Value2 = 1
$ws.Cells.Item(50, 18).Value2 = "Hortaliza"

# Row 51
$ws.Cells.Item(51, 1).Value2 = 1
$ws.Cells.Item(51, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(51, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(51, 4).Value2 = 44412
$ws.Cells.Item(51, 5).Value2 = 15
$ws.Cells.Item(51, 6).Value2 = 100112031
$ws.Cells.Item(51, 7).Value2 = "Poroto verde"
$ws.Cells.Item(51, 8).Value2 = "Magnum"
$ws.Cells.Item(51, 9).Value2 = "Primera"
$ws.Cells.Item(51, 10).Value2 = 250
$ws.Cells.Item(51, 11).Value2 = 28000
$ws.Cells.Item(51, 12).Value2 = 29000
$ws.Cells.Item(51, 13).Value2 = 28500
$ws.Cells.Item(51, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(51, 15).Value2 = "Perú"
$ws.Cells.Item(51, 16).Value2 = 1140
$ws.Cells.Item(51, 17).Value2 = 25
$ws.Cells.Item(51, 18).Value2 = "Hortaliza"

# Row 52
$ws.Cells.Item(52, 1).Value2 = 1
$ws.Cells.Item(52, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(52, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(52, 4).Value2 = 44274
$ws.Cells.Item(52, 5).Value2 = 15
$ws.Cells.Item(52, 6).Value2 = 100112031
$ws.Cells.Item(52, 7).Value2 = "Poroto verde"
$ws.Cells.Item(52, 8).Value2 = "Sin especificar"
$ws.Cells.Item(52, 9).Value2 = "Primera"
$ws.Cells.Item(52, 10).Value2 = 1500
$ws.Cells.Item(52, 11).Value2 = 2200
$ws.Cells.Item(52, 12).Value2 = 2300
$ws.Cells.Item(52, 13).Value2 = 2250
$ws.Cells.Item(52, 14).Value2 = "$/kilo"
$ws.Cells.Item(52, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(52, 16).Value2 = 2250
$ws.Cells.Item(52, 17).Value2 = 1
$ws.Cells.Item(52, 18).Value2 = "Hortaliza"

# Row 53
$ws.Cells.Item(53, 1).Value2 = 1
$ws.Cells.Item(53, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(53, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(53, 4).Value2 = 44483
$ws.Cells.Item(53, 5).Value2 = 15
$ws.Cells.Item(53, 6).Value2 = 100112031
$ws.Cells.Item(53, 7).Value2 = "Poroto verde"
$ws.Cells.Item(53, 8).Value2 = "Sin especificar"
$ws.Cells.Item(53, 9).Value2 = "Primera"
$ws.Cells.Item(53, 10).Value2 = 1100
$ws.Cells.Item(53, 11).Value2 = 1200
$ws.Cells.Item(53, 12).Value2 = 1300
$ws.Cells.Item(53, 13).Value2 = 1250
$ws.Cells.Item(53, 14).Value2 = "$/kilo"
$ws.Cells.Item(53, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(53, 16).Value2 = 1250
$ws.Cells.Item(53, 17).Value2 = 1
$ws.Cells.Item(53, 18).Value2 = "Hortaliza"

# Row 54
$ws.Cells.Item(54, 1).Value2 = 1
$ws.Cells.Item(54, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(54, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(54, 4).Value2 = 44371
$ws.Cells.Item(54, 5).Value2 = 15
$ws.Cells.Item(54, 6).Value2 = 100112031
$ws.Cells.Item(54, 7).Value2 = "Poroto verde"
$ws.Cells.Item(54, 8).Value2 = "Sin especificar"
$ws.Cells.Item(54, 9).Value2 = "Primera"
$ws.Cells.Item(54, 10).Value2 = 1200
$ws.Cells.Item(54, 11).Value2 = 900
$ws.Cells.Item(54, 12).Value2 = 1000
$ws.Cells.Item(54, 13).Value2 = 942
$ws.Cells.Item(54, 14).Value2 = "$/kilo"
$ws.Cells.Item(54, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(54, 16).Value2 = 942
$ws.Cells.Item(54, 17).Value2 = 1
$ws.Cells.Item(54, 18).Value2 = "Hortaliza"

# Row 55
$ws.Cells.Item(55, 1).Value2 = 1
$ws.Cells.Item(55, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(55, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(55, 4).Value2 = 44305
$ws.Cells.Item(55, 5).Value2 = 15
$ws.Cells.Item(55, 6).Value2 = 100112031
$ws.Cells.Item(55, 7).Value2 = "Poroto verde"
$ws.Cells.Item(55, 8).Value2 = "Sin especificar"
$ws.Cells.Item(55, 9).Value2 = "Primera"
$ws.Cells.Item(55, 10).Value2 = 1300
$ws.Cells.Item(55, 11).Value2 = 850
$ws.Cells.Item(55, 12).Value2 = 900
$ws.Cells.Item(55, 13).Value2 = 875
$ws.Cells.Item(55, 14).Value2 = "$/kilo"
$ws.Cells.Item(55, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(55, 16).Value2 = 875
$ws.Cells.Item(55, 17).Value2 = 1
$ws.Cells.Item(55, 18).Value2 = "Hortaliza"

# Row 56
$ws.Cells.Item(56, 1).Value2 = 1
$ws.Cells.Item(56, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(56, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(56, 4).Value2 = 44385
$ws.Cells.Item(56, 5).Value2 = 15
$ws.Cells.Item(56, 6).Value2 = 100112031
$ws.Cells.Item(56, 7).Value2 = "Poroto verde"
$ws.Cells.Item(56, 8).Value2 = "Sin especificar"
$ws.Cells.Item(56, 9).Value2 = "Primera"
$ws.Cells.Item(56, 10).Value2 = 1500
$ws.Cells.Item(56, 11).Value2 = 1200
$ws.Cells.Item(56, 12).Value2 = 1300
$ws.Cells.Item(56, 13).Value2 = 1250
$ws.Cells.Item(56, 14).Value2 = "$/kilo"
$ws.Cells.Item(56, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(56, 16).Value2 = 1250
$ws.Cells.Item(56, 17).Value2 = 1
$ws.Cells.Item(56, 18).Value2 = "Hortaliza"

# Row 57
$ws.Cells.Item(57, 1).Value2 = 1
$ws.Cells.Item(57, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(57, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(57, 4).Value2 = 44641
$ws.Cells.Item(57, 5).Value2 = 15
$ws.Cells.Item(57, 6).Value2 = 100112031
$ws.Cells.Item(57, 7).Value2 = "Poroto verde"
$ws.Cells.Item(57, 8).Value2 = "Sin especificar"
$ws.Cells.Item(57, 9).Value2 = "Primera"
$ws.Cells.Item(57, 10).Value2 = 1000
$ws.Cells.Item(57, 11).Value2 = 1300
$ws.Cells.Item(57, 12).Value2 = 1400
$ws.Cells.Item(57, 13).Value2 = 1350
$ws.Cells.Item(57, 14).Value2 = "$/kilo"
$ws.Cells.Item(57, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(57, 16).Value2 = 1350
$ws.Cells.Item(57, 17).Value2 = 1
$ws.Cells.Item(57, 18).Value2 = "Hortaliza"

# Row 58
$ws.Cells.Item(58, 1).Value2 = 1
$ws.Cells.Item(58, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(58, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(58, 4).Value2 = 44323
$ws.Cells.Item(58, 5).Value2 = 15
$ws.Cells.Item(58, 6).Value2 = 100112031
$ws.Cells.Item(58, 7).Value2 = "Poroto verde"
$ws.Cells.Item(58, 8).Value2 = "Sin especificar"
$ws.Cells.Item(58, 9).Value2 = "Primera"
$ws.Cells.Item(58, 10).Value2 = 1500
$ws.Cells.Item(58, 11).Value2 = 1700
$ws.Cells.Item(58, 12).Value2 = 1800
$ws.Cells.Item(58, 13).Value2 = 1750
$ws.Cells.Item(58, 14).Value2 = "$/kilo"
$ws.Cells.Item(58, 15).Value2 = "Perú"
$ws.Cells.Item(58, 16).Value2 = 1750
$ws.Cells.Item(58, 17).Value2 = 1
$ws.Cells.Item(58, 18).Value2 = "Hortaliza"

# Row 59
$ws.Cells.Item(59, 1).Value2 = 1
$ws.Cells.Item(59, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(59, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(59, 4).Value2 = 44889
$ws.Cells.Item(59, 5).Value2 = 15
$ws.Cells.Item(59, 6).Value2 = 100112031
$ws.Cells.Item(59, 7).Value2 = "Poroto verde"
$ws.Cells.Item(59, 8).Value2 = "Sin especificar"
$ws.Cells.Item(59, 9).Value2 = "Primera"
$ws.Cells.Item(59, 10).Value2 = 500
$ws.Cells.Item(59, 11).Value2 = 900
$ws.Cells.Item(59, 12).Value2 = 1000
$ws.Cells.Item(59, 13).Value2 = 950
$ws.Cells.Item(59, 14).Value2 = "$/kilo"
$ws.Cells.Item(59, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(59, 16).Value2 = 950
$ws.Cells.Item(59, 17).Value2 = 1
$ws.Cells.Item(59, 18).Value2 = "Hortaliza"

# Row 60
$ws.Cells.Item(60, 1).Value2 = 1
$ws.Cells.Item(60, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(60, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(60, 4).Value2 = 44649
$ws.Cells.Item(60, 5).Value2 = 15
$ws.Cells.Item(60, 6).Value2 = 100112031
$ws.Cells.Item(60, 7).Value2 = "Poroto verde"
$ws.Cells.Item(60, 8).Value2 = "Sin especificar"
$ws.Cells.Item(60, 9).Value2 = "Primera"
$ws.Cells.Item(60, 10).Value2 = 1700
$ws.Cells.Item(60, 11).Value2 = 500
$ws.Cells.Item(60, 12).Value2 = 600
$ws.Cells.Item(60, 13).Value2 = 550
$ws.Cells.Item(60, 14).Value2 = "$/kilo"
$ws.Cells.Item(60, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 16).Value2 = 550
$ws.Cells.Item(60, 17).Value2 = 1
$ws.Cells.Item(60, 18).Value2 = "Hortaliza"

# Row 61
$ws.Cells.Item(61, 1).Value2 = 1
$ws.Cells.Item(61, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(61, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(61, 4).Value2 = 44704
$ws.Cells.Item(61, 5).Value2 = 15
$ws.Cells.Item(61, 6).Value2 = 100112031
$ws.Cells.Item(61, 7).Value2 = "Poroto verde"
$ws.Cells.Item(61, 8).Value2 = "Magnum"
$ws.Cells.Item(61, 9).Value2 = "Primera"
$ws.Cells.Item(61, 10).Value2 = 200
$ws.Cells.Item(61, 11).Value2 = 19000
$ws.Cells.Item(61, 12).Value2 = 20000
$ws.Cells.Item(61, 13).Value2 = 19500
$ws.Cells.Item(61, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(61, 15).Value2 = "Perú"
$ws.Cells.Item(61, 16).Value2 = 780
$ws.Cells.Item(61, 17).Value2 = 25
$ws.Cells.Item(61, 18).Value2 = "Hortaliza"

# Row 62
$ws.Cells.Item(62, 1).Value2 = 1
$ws.Cells.Item(62, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(62, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(62, 4).Value2 = 44370
$ws.Cells.Item(62, 5).Value2 = 15
$ws.Cells.Item(62, 6).Value2 = 100112031
$ws.Cells.Item(62, 7).Value2 = "Poroto verde"
$ws.Cells.Item(62, 8).Value2 = "Magnum"
$ws.Cells.Item(62, 9).Value2 = "Primera"
$ws.Cells.Item(62, 10).Value2 = 80
$ws.Cells.Item(62, 11).Value2 = 19000
$ws.Cells.Item(62, 12).Value2 = 20000
$ws.Cells.Item(62, 13).Value2 = 19375
$ws.Cells.Item(62, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(62, 15).Value2 = "Perú"
$ws.Cells.Item(62, 16).Value2 = 775
$ws.Cells.Item(62, 17).Value2 = 25
$ws.Cells.Item(62, 18).Value2 = "Hortaliza"

# Row 63
$ws.Cells.Item(63, 1).Value2 = 1
$ws.Cells.Item(63, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(63, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(63, 4).Value2 = 44370
$ws.Cells.Item(63, 5).Value2 = 15
$ws.Cells.Item(63, 6).Value2 = 100112031
$ws.Cells.Item(63, 7).Value2 = "Poroto verde"
$ws.Cells.Item(63, 8).Value2 = "Magnum"
$ws.Cells.Item(63, 9).Value2 = "Segunda"
$ws.Cells.Item(63, 10).Value2 = 40
$ws.Cells.Item(63, 11).Value2 = 17000
$ws.Cells.Item(63, 12).Value2 = 18000
$ws.Cells.Item(63, 13).Value2 = 17500
$ws.Cells.Item(63, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(63, 15).Value2 = "Perú"
$ws.Cells.Item(63, 16).Value2 = 700
$ws.Cells.Item(63, 17).Value2 = 25
$ws.Cells.Item(63, 18).Value2 = "Hortaliza"

# Row 64
$ws.Cells.Item(64, 1).Value2 = 1
$ws.Cells.Item(64, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(64, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(64, 4).Value2 = 44334
$ws.Cells.Item(64, 5).Value2 = 15
$ws.Cells.Item(64, 6).Value2 = 100112031
$ws.Cells.Item(64, 7).Value2 = "Poroto verde"
$ws.Cells.Item(64, 8).Value2 = "Sin especificar"
$ws.Cells.Item(64, 9).Value2 = "Primera"
$ws.Cells.Item(64, 10).Value2 = 1440
$ws.Cells.Item(64, 11).Value2 = 1300
$ws.Cells.Item(64, 12).Value2 = 1400
$ws.Cells.Item(64, 13).Value2 = 1350
$ws.Cells.Item(64, 14).Value2 = "$/kilo"
$ws.Cells.Item(64, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(64, 16).Value2 = 1350
$ws.Cells.Item(64, 17).Value2 = 1
$ws.Cells.Item(64, 18).Value2 = "Hortaliza"

# Row 65
$ws.Cells.Item(65, 1).Value2 = 1
$ws.Cells.Item(65, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(65, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(65, 4).Value2 = 44550
$ws.Cells.Item(65, 5).Value2 = 15
$ws.Cells.Item(65, 6).Value2 = 100112031
$ws.Cells.Item(65, 7).Value2 = "Poroto verde"
$ws.Cells.Item(65, 8).Value2 = "Sin especificar"
$ws.Cells.Item(65, 9).Value2 = "Primera"
$ws.Cells.Item(65, 10).Value2 = 1500
$ws.Cells.Item(65, 11).Value2 = 300
$ws.Cells.Item(65, 12).Value2 = 350
$ws.Cells.Item(65, 13).Value2 = 325
$ws.Cells.Item(65, 14).Value2 = "$/kilo"
$ws.Cells.Item(65, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(65, 16).Value2 = 325
$ws.Cells.Item(65, 17).Value2 = 1
$ws.Cells.Item(65, 18).Value2 = "Hortaliza"

# Row 66
$ws.Cells.Item(66, 1).Value2 = 1
$ws.Cells.Item(66, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(66, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(66, 4).Value2 = 44204
$ws.Cells.Item(66, 5).Value2 = 15
$ws.Cells.Item(66, 6).Value2 = 100112031
$ws.Cells.Item(66, 7).Value2 = "Poroto verde"
$ws.Cells.Item(66, 8).Value2 = "Sin especificar"
$ws.Cells.Item(66, 9).Value2 = "Primera"
$ws.Cells.Item(66, 10).Value2 = 1500
$ws.Cells.Item(66, 11).Value2 = 1000
$ws.Cells.Item(66, 12).Value2 = 1100
$ws.Cells.Item(66, 13).Value2 = 1050
$ws.Cells.Item(66, 14).Value2 = "$/kilo"
$ws.Cells.Item(66, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(66, 16).Value2 = 1050
$ws.Cells.Item(66, 17).Value2 = 1
$ws.Cells.Item(66, 18).Value2 = "Hortaliza"

# Row 67
$ws.Cells.Item(67, 1).Value2 = 1
$ws.Cells.Item(67, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(67, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(67, 4).Value2 = 44656
$ws.Cells.Item(67, 5).Value2 = 15
$ws.Cells.Item(67, 6).Value2 = 100112031
$ws.Cells.Item(67, 7).Value2 = "Poroto verde"
$ws.Cells.Item(67, 8).Value2 = "Sin especificar"
$ws.Cells.Item(67, 9).Value2 = "Primera"
$ws.Cells.Item(67, 10).Value2 = 1500
$ws.Cells.Item(67, 11).Value2 = 600
$ws.Cells.Item(67, 12).Value2 = 700
$ws.Cells.Item(67, 13).Value2 = 650
$ws.Cells.Item(67, 14).Value2 = "$/kilo"
$ws.Cells.Item(67, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(67, 16).Value2 = 650
$ws.Cells.Item(67, 17).Value2 = 1
$ws.Cells.Item(67, 18).Value2 = "Hortaliza"

# Row 68
$ws.Cells.Item(68, 1).Value2 = 1
$ws.Cells.Item(68, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(68, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(68, 4).Value2 = 44278
$ws.Cells.Item(68, 5).Value2 = 15
$ws.Cells.Item(68, 6).Value2 = 100112031
$ws.Cells.Item(68, 7).Value2 = "Poroto verde"
$ws.Cells.Item(68, 8).Value2 = "Sin especificar"
$ws.Cells.Item(68, 9).Value2 = "Primera"
$ws.Cells.Item(68, 10).Value2 = 1300
$ws.Cells.Item(68, 11).Value2 = 2400
$ws.Cells.Item(68, 12).Value2 = 2500
$ws.Cells.Item(68, 13).Value2 = 2450
$ws.Cells.Item(68, 14).Value2 = "$/kilo"
$ws.Cells.Item(68, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value2 = 2450
$ws.Cells.Item(68, 17).Value2 = 1
$ws.Cells.Item(68, 18).Value2 = "Hortaliza"

# Row 69
$ws.Cells.Item(69, 1).Value2 = 1
$ws.Cells.Item(69, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(69, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(69, 4).Value2 = 44298
$ws.Cells.Item(69, 5).Value2 = 15
$ws.Cells.Item(69, 6).Value2 = 100112031
$ws.Cells.Item(69, 7).Value2 = "Poroto verde"
$ws.Cells.Item(69, 8).Value2 = "Sin especificar"
$ws.Cells.Item(69, 9).Value2 = "Primera"
$ws.Cells.Item(69, 10).Value2 = 1500
$ws.Cells.Item(69, 11).Value2 = 700
$ws.Cells.Item(69, 12).Value2 = 800
$ws.Cells.Item(69, 13).Value2 = 750
$ws.Cells.Item(69, 14).Value2 = "$/kilo"
$ws.Cells.Item(69, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(69, 16).Value2 = 750
$ws.Cells.Item(69, 17).Value2 = 1
$ws.Cells.Item(69, 18).Value2 = "Hortaliza"

# Row 70
$ws.Cells.Item(70, 1).Value2 = 1
$ws.Cells.Item(70, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(70, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(70, 4).Value2 = 44162
$ws.Cells.Item(70, 5).Value2 = 15
$ws.Cells.Item(70, 6).Value2 = 100112031
$ws.Cells.Item(70, 7).Value2 = "Poroto verde"
$ws.Cells.Item(70, 8).Value2 = "Sin especificar"
$ws.Cells.Item(70, 9).Value2 = "Primera"
$ws.Cells.Item(70, 10).Value2 = 1870
$ws.Cells.Item(70, 11).Value2 = 950
$ws.Cells.Item(70, 12).Value2 = 1000
$ws.Cells.Item(70, 13).Value2 = 975
$ws.Cells.Item(70, 14).Value2 = "$/kilo"
$ws.Cells.Item(70, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(70, 16).Value2 = 975
$ws.Cells.Item(70, 17).Value2 = 1
$ws.Cells.Item(70, 18).Value2 = "Hortaliza"

# Row 71
$ws.Cells.Item(71, 1).Value2 = 1
$ws.Cells.Item(71, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(71, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(71, 4).Value2 = 44291
$ws.Cells.Item(71, 5).Value2 = 15
$ws.Cells.Item(71, 6).Value2 = 100112031
$ws.Cells.Item(71, 7).Value2 = "Poroto verde"
$ws.Cells.Item(71, 8).Value2 = "Sin especificar"
$ws.Cells.Item(71, 9).Value2 = "Primera"
$ws.Cells.Item(71, 10).Value2 = 1500
$ws.Cells.Item(71, 11).Value2 = 1400
$ws.Cells.Item(71, 12).Value2 = 1500
$ws.Cells.Item(71, 13).Value2 = 1450
$ws.Cells.Item(71, 14).Value2 = "$/kilo"
$ws.Cells.Item(71, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value2 = 1450
$ws.Cells.Item(71, 17).Value2 = 1
$ws.Cells.Item(71, 18).Value2 = "Hortaliza"

# Row 72
$ws.Cells.Item(72, 1).Value2 = 1
$ws.Cells.Item(72, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value2 = 44778
$ws.Cells.Item(72, 5).Value2 = 15
$ws.Cells.Item(72, 6).Value2 = 100112031
$ws.Cells.Item(72, 7).Value2 = "Poroto verde"
$ws.Cells.Item(72, 8).Value2 = "Magnum"
$ws.Cells.Item(72, 9).Value2 = "Primera"
$ws.Cells.Item(72, 10).Value2 = 300
$ws.Cells.Item(72, 11).Value2 = 27000
$ws.Cells.Item(72, 12).Value2 = 28000
$ws.Cells.Item(72, 13).Value2 = 27500
$ws.Cells.Item(72, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(72, 15).Value2 = "Perú"
$ws.Cells.Item(72, 16).Value2 = 1100
$ws.Cells.Item(72, 17).Value2 = 25
$ws.Cells.Item(72, 18).Value2 = "Hortaliza"

# Row 73
$ws.Cells.Item(73, 1).Value2 = 1
$ws.Cells.Item(73, 2).Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(73, 3).Value2 = "Arica y Parinacota"
$ws.Cells.Item(73, 4).Value2 = 44496
$ws.Cells.Item(73, 5).Value2 = 15
$ws.Cells.Item(73, 6).Value2 = 100112031
$ws.Cells.Item(73, 7).Value2 = "Poroto verde"
$ws.Cells.Item(73, 8).Value2 = "Sin especificar"
$ws.Cells.Item(73, 9).Value2 = "Primera"
$ws.Cells.Item(73, 10).Value2 = 1300
$ws.Cells.Item(73, 11).Value2 = 1400
$ws.Cells.Item(73, 12).Value2 = 1500
$ws.Cells.Item(73, 13).Value2 = 1450
$ws.Cells.Item(73, 14).Value2 = "$/kilo"
$ws.Cells.Item(73, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(73, 16).Value2 = 1450
$ws.Cells.Item(73, 17).Value2 = 1
$ws.Cells.Item(73, 18).Value2 = "Hortaliza"
